$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("for paper")
$cell = $ws.Range("N50")
Write-Host "Formula:" $cell.Formula
Write-Host "FormulaR1C1:" $cell.FormulaR1C1
